$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2.433531715253719
